$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "Producto" column (column B) entirely -- Material/Cantidad shift left
$ws.Range("B1").EntireColumn.Delete()

# Delete the now-empty rows below the header row (rows 2 through 14)
$ws.Range("A2:A14").EntireRow.Delete()

# Set the active selection to C9 as shown in the diff
$ws.Range("C9").Select()
